$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell updates per the source diff. Cells whose new text looks like a
# plain number (e.g. "319.70", "1.001") are forced to Text format first,
# so Excel keeps them as literal strings (matching the original text
# storage) instead of silently coercing them into floating point
# numbers and losing formatting like trailing zeros.

$ws.Range('D2').Value = '29.928.11'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '1.909.28'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '319.70'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5049'
$ws.Range('E7').Value = '  -2.32%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4053'
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08272'
$ws.Range('E9').Value = '  -1.96%  '
$ws.Range('E10').Value = '  -1.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.92'
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('E12').Value = '  +4.01%  '
$ws.Range('D13').Value = '1.905.70'
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.397'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.206'
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.9987'
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '92.09'
$ws.Range('E17').Value = '  -2.76%  '
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06506'
$ws.Range('E19').Value = '  -2.42%  '
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.933'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').Value = '29.969.82'
$ws.Range('E23').Value = '  -0.74%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.30'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.200'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '22.37'
$ws.Range('E26').Value = '  +2.84%  '
$ws.Range('D27').Value = '2.126.60'
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '162.11'
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('E29').Value = '  -3.62%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '128.92'
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1036'
$ws.Range('E32').Value = '  -1.88%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.926'
$ws.Range('E33').Value = '  -1.94%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.813'
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.405'
$ws.Range('E35').Value = '  +2.75%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02440'
$ws.Range('E36').Value = '  -1.85%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06401'
$ws.Range('E37').Value = '  -2.36%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2143'
$ws.Range('E38').Value = '  -2.88%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.708'
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.193'
$ws.Range('E40').Value = '  -1.95%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6467'
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('E42').Value = '  -3.40%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.209'
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.222'
$ws.Range('E44').Value = '  +8.11%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.24'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6030'
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.635'
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '122.23'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.204'
$ws.Range('E49').Value = '  -2.37%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '78.69'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.125'
$ws.Range('E51').Value = '  -3.46%  '
